$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D ("Product" header / {booking:product_or_services} data)
# This shifts old column E -> D and old column F -> E
$ws.Columns("D:D").Delete()
